# "Fixed HKCD excel file"
#
# The sheet was originally authored/saved from a Chinese-locale Excel build
# (tab name "工作表1", selection left on A183). This pass renames the sheet
# back to the generic "Sheet1" and resets the saved cursor position, matching
# the state captured after the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "工作表1" -> "Sheet1"
$ws.Name = "Sheet1"

# Move/record the selected cell as R18 (was A183)
[void]$ws.Range("R18").Select()
